$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = [double]"22.88000000000014"
$ws.Range("H2").Value = [double]"0.002969250382238275"
$ws.Range("I2").Value = [double]"0.002969250382238275"
$ws.Range("L2").Value = [double]"29.12558995404007"
$ws.Range("M2").Value = "[8.442683135941081, 49.80849677213906]"
$ws.Range("N2").Value = [double]"0.006817996700800411"
$ws.Range("O2").Value = [double]"0.006817996700800411"
$ws.Range("P2").Value = [double]"1.930868758062195"
$ws.Range("Q2").Value = "[1.1383949355350396, 2.7233425805893496]"
$ws.Range("R2").Value = [double]"1.251423347237157e-05"
$ws.Range("S2").Value = [double]"1.251423347237157e-05"
$ws.Range("T2").Value = [double]"56.81737397420945"
$ws.Range("U2").Value = "[45.58417061611381, 68.05057733230508]"
$ws.Range("V2").Value = [double]"2.897682094271659e-13"
$ws.Range("W2").Value = [double]"2.897682094271659e-13"
$ws.Range("X2").Value = [double]"15.84880880880891"
$ws.Range("Y2").Value = [double]"12.96304304304312"
$ws.Range("Z2").Value = [double]"18.73457457457469"
# Row 3
$ws.Range("F3").Value = [double]"22.88000000000014"
$ws.Range("H3").Value = [double]"0.001435797550979778"
$ws.Range("I3").Value = [double]"0.001435797550979778"
$ws.Range("L3").Value = [double]"42.10747065432017"
$ws.Range("M3").Value = "[12.735857552629938, 71.47908375601041]"
$ws.Range("N3").Value = [double]"0.005949582198544023"
$ws.Range("O3").Value = [double]"0.005949582198544023"
$ws.Range("P3").Value = [double]"2.044079304137503"
$ws.Range("Q3").Value = "[1.3270791789938858, 2.7610794292811196]"
$ws.Range("R3").Value = [double]"7.531544425898318e-07"
$ws.Range("S3").Value = [double]"7.531544425898318e-07"
$ws.Range("T3").Value = [double]"54.27314089426032"
$ws.Range("U3").Value = "[39.02388308216193, 69.5223987063587]"
$ws.Range("V3").Value = [double]"5.738788999565259e-09"
$ws.Range("W3").Value = [double]"5.738788999565259e-09"
$ws.Range("X3").Value = [double]"15.43655655655665"
$ws.Range("Y3").Value = [double]"12.8256256256257"
$ws.Range("Z3").Value = [double]"18.0474874874876"
# Row 4
$ws.Range("F4").Value = [double]"22.88000000000014"
$ws.Range("H4").Value = [double]"1.692863753621765e-07"
$ws.Range("I4").Value = [double]"1.692863753621765e-07"
$ws.Range("L4").Value = [double]"47.73523877690788"
$ws.Range("M4").Value = "[28.075275961760816, 67.39520159205495]"
$ws.Range("N4").Value = [double]"1.324170247030132e-05"
$ws.Range("O4").Value = [double]"1.324170247030132e-05"
$ws.Range("P4").Value = [double]"2.169868799776734"
$ws.Range("Q4").Value = "[1.7673424137311953, 2.5723951858222724]"
$ws.Range("R4").Value = [double]"3.708144902248023e-14"
$ws.Range("S4").Value = [double]"3.708144902248023e-14"
$ws.Range("T4").Value = [double]"56.81463105031748"
$ws.Range("U4").Value = "[46.71769237313327, 66.9115697275017]"
$ws.Range("V4").Value = [double]"8.881784197001252e-15"
$ws.Range("W4").Value = [double]"8.881784197001252e-15"
$ws.Range("X4").Value = [double]"14.97849849849859"
$ws.Range("Y4").Value = [double]"13.5127127127128"
$ws.Range("Z4").Value = [double]"16.44428428428438"
# Row 5
$ws.Range("F5").Value = [double]"22.88000000000014"
$ws.Range("H5").Value = [double]"0.001191809781018405"
$ws.Range("I5").Value = [double]"0.001191809781018405"
$ws.Range("L5").Value = [double]"40.47037597884702"
$ws.Range("M5").Value = "[13.28984726424791, 67.65090469344612]"
$ws.Range("N5").Value = [double]"0.004402817195519271"
$ws.Range("O5").Value = [double]"0.004402817195519271"
$ws.Range("P5").Value = [double]"2.698184681461504"
$ws.Range("Q5").Value = "[1.9937635058818106, 3.4026058570411974]"
$ws.Range("R5").Value = [double]"8.995379996434849e-10"
$ws.Range("S5").Value = [double]"8.995379996434849e-10"
$ws.Range("T5").Value = [double]"60.66140864356841"
$ws.Range("U5").Value = "[45.864693202043405, 75.45812408509343]"
$ws.Range("V5").Value = [double]"1.460569443167969e-10"
$ws.Range("W5").Value = [double]"1.460569443167969e-10"
$ws.Range("X5").Value = [double]"13.05465465465473"
$ws.Range("Y5").Value = [double]"10.48952952952959"
$ws.Range("Z5").Value = [double]"15.61977977977987"
# Row 6
$ws.Range("F6").Value = [double]"22.88000000000014"
$ws.Range("H6").Value = [double]"0.004422092388514876"
$ws.Range("I6").Value = [double]"0.004422092388514876"
$ws.Range("L6").Value = [double]"30.67493491087199"
$ws.Range("M6").Value = "[8.92240757082233, 52.42746225092166]"
$ws.Range("N6").Value = [double]"0.006746261408696164"
$ws.Range("O6").Value = [double]"0.006746261408696164"
$ws.Range("P6").Value = [double]"3.025237370123504"
$ws.Range("Q6").Value = "[2.0943951023931953, 3.9560796378538137]"
$ws.Range("R6").Value = [double]"4.818426413422117e-08"
$ws.Range("S6").Value = [double]"4.818426413422117e-08"
$ws.Range("T6").Value = [double]"58.11403960619087"
$ws.Range("U6").Value = "[45.35965897894101, 70.86842023344073]"
$ws.Range("V6").Value = [double]"7.137623825315131e-12"
$ws.Range("W6").Value = [double]"7.137623825315131e-12"
$ws.Range("X6").Value = [double]"11.86370370370377"
$ws.Range("Y6").Value = [double]"8.474074074074123"
$ws.Range("Z6").Value = [double]"15.25333333333342"
# Row 7
$ws.Range("F7").Value = [double]"22.88000000000014"
$ws.Range("H7").Value = [double]"3.663884824522956e-05"
$ws.Range("I7").Value = [double]"3.663884824522956e-05"
$ws.Range("L7").Value = [double]"49.54763798415998"
$ws.Range("M7").Value = "[26.708251855548838, 72.38702411277113]"
$ws.Range("N7").Value = [double]"7.253278232366434e-05"
$ws.Range("O7").Value = [double]"7.253278232366434e-05"
$ws.Range("P7").Value = [double]"-2.956053147521927"
$ws.Range("Q7").Value = "[-3.5472637770263122, -2.3648425180175425]"
$ws.Range("R7").Value = [double]"4.172218126541338e-13"
$ws.Range("S7").Value = [double]"4.172218126541338e-13"
$ws.Range("T7").Value = [double]"61.07628863656355"
$ws.Range("U7").Value = "[46.94132900062601, 75.2112482725011]"
$ws.Range("V7").Value = [double]"3.346278809601699e-11"
$ws.Range("W7").Value = [double]"3.346278809601699e-11"
$ws.Range("X7").Value = [double]"10.76436436436443"
$ws.Range("Y7").Value = [double]"8.611491491491545"
$ws.Range("Z7").Value = [double]"12.91723723723731"
# Row 8
$ws.Range("F8").Value = [double]"25.8500000000006"
$ws.Range("H8").Value = [double]"2.197977608209278e-05"
$ws.Range("I8").Value = [double]"2.197977608209278e-05"
$ws.Range("L8").Value = [double]"48.43851995883746"
$ws.Range("M8").Value = "[23.150702177709107, 73.72633773996581]"
$ws.Range("N8").Value = [double]"0.0003617639804807116"
$ws.Range("O8").Value = [double]"0.0003617639804807116"
$ws.Range("P8").Value = [double]"-2.528368862348542"
$ws.Range("Q8").Value = "[-3.0315268449054646, -2.025210879791619]"
$ws.Range("R8").Value = [double]"3.563815909046752e-13"
$ws.Range("S8").Value = [double]"3.563815909046752e-13"
$ws.Range("T8").Value = [double]"55.54862666724694"
$ws.Range("U8").Value = "[42.796984258502974, 68.30026907599091]"
$ws.Range("V8").Value = [double]"2.651012742660441e-11"
$ws.Range("W8").Value = [double]"2.651012742660441e-11"
$ws.Range("X8").Value = [double]"10.40210210210234"
$ws.Range("Y8").Value = [double]"8.332032032032229"
$ws.Range("Z8").Value = [double]"12.47217217217246"
# Row 9
$ws.Range("F9").Value = [double]"25.8500000000006"
$ws.Range("H9").Value = [double]"8.805004662726379e-05"
$ws.Range("I9").Value = [double]"8.805004662726379e-05"
$ws.Range("L9").Value = [double]"53.80966148241092"
$ws.Range("M9").Value = "[27.05353053453014, 80.5657924302917]"
$ws.Range("N9").Value = [double]"0.0001992175710661126"
$ws.Range("O9").Value = [double]"0.0001992175710661126"
$ws.Range("P9").Value = [double]"-2.993789996213697"
$ws.Range("Q9").Value = "[-3.5724216761541596, -2.415158316273234]"
$ws.Range("R9").Value = [double]"1.407762795224698e-13"
$ws.Range("S9").Value = [double]"1.407762795224698e-13"
$ws.Range("T9").Value = [double]"57.63289978113562"
$ws.Range("U9").Value = "[42.28788658772592, 72.97791297454532]"
$ws.Range("V9").Value = [double]"1.494030454907147e-09"
$ws.Range("W9").Value = [double]"1.494030454907147e-09"
$ws.Range("X9").Value = [double]"12.3169169169172"
$ws.Range("Y9").Value = [double]"9.936336336336566"
$ws.Range("Z9").Value = [double]"14.69749749749784"
# Row 10
$ws.Range("F10").Value = [double]"25.8500000000006"
$ws.Range("H10").Value = [double]"9.57815707210985e-05"
$ws.Range("I10").Value = [double]"9.57815707210985e-05"
$ws.Range("L10").Value = [double]"49.55127281344329"
$ws.Range("M10").Value = "[25.21324525601446, 73.88930037087212]"
$ws.Range("N10").Value = [double]"0.0001703060061821748"
$ws.Range("O10").Value = [double]"0.0001703060061821748"
$ws.Range("P10").Value = [double]"2.987500521431735"
$ws.Range("Q10").Value = "[2.408868841491273, 3.5661322013721968]"
$ws.Range("R10").Value = [double]"1.505462421391712e-13"
$ws.Range("S10").Value = [double]"1.505462421391712e-13"
$ws.Range("T10").Value = [double]"63.53142108197113"
$ws.Range("U10").Value = "[49.13502097104352, 77.92782119289875]"
$ws.Range("V10").Value = [double]"1.823519113486327e-11"
$ws.Range("W10").Value = [double]"1.823519113486327e-11"
$ws.Range("X10").Value = [double]"13.55895895895928"
$ws.Range("Y10").Value = [double]"11.17837837837864"
$ws.Range("Z10").Value = [double]"15.93953953953991"
# Row 11
$ws.Range("F11").Value = [double]"25.8500000000006"
$ws.Range("H11").Value = [double]"5.030264001670837e-06"
$ws.Range("I11").Value = [double]"5.030264001670837e-06"
$ws.Range("L11").Value = [double]"48.44643120871829"
$ws.Range("M11").Value = "[23.994541340345602, 72.89832107709098]"
$ws.Range("N11").Value = [double]"0.0002402387540294537"
$ws.Range("O11").Value = [double]"0.0002402387540294537"
$ws.Range("P11").Value = [double]"2.371131992799503"
$ws.Range("Q11").Value = "[1.8931319093704255, 2.849132076228581]"
$ws.Range("R11").Value = [double]"5.349054532644004e-13"
$ws.Range("S11").Value = [double]"5.349054532644004e-13"
$ws.Range("T11").Value = [double]"64.77543643679314"
$ws.Range("U11").Value = "[52.37722423385698, 77.17364863972931]"
$ws.Range("V11").Value = [double]"1.028066520802895e-13"
$ws.Range("W11").Value = [double]"1.028066520802895e-13"
$ws.Range("X11").Value = [double]"16.09479479479517"
$ws.Range("Y11").Value = [double]"14.12822822822856"
$ws.Range("Z11").Value = [double]"18.06136136136178"
# Row 12
$ws.Range("F12").Value = [double]"25.8500000000006"
$ws.Range("H12").Value = [double]"6.427239626405346e-06"
$ws.Range("I12").Value = [double]"6.427239626405346e-06"
$ws.Range("L12").Value = [double]"49.12674119762966"
$ws.Range("M12").Value = "[27.20322537690768, 71.05025701835164]"
$ws.Range("N12").Value = [double]"4.559940512183225e-05"
$ws.Range("O12").Value = [double]"4.559940512183225e-05"
$ws.Range("P12").Value = [double]"1.767342413731195"
$ws.Range("Q12").Value = "[1.2390265320464255, 2.295658295415965]"
$ws.Range("R12").Value = [double]"2.49845932831505e-08"
$ws.Range("S12").Value = [double]"2.49845932831505e-08"
$ws.Range("T12").Value = [double]"63.56828342576374"
$ws.Range("U12").Value = "[50.58624727703989, 76.55031957448759]"
$ws.Range("V12").Value = [double]"8.01358979174438e-13"
$ws.Range("W12").Value = [double]"8.01358979174438e-13"
$ws.Range("X12").Value = [double]"18.57887887887931"
$ws.Range("Y12").Value = [double]"16.40530530530568"
$ws.Range("Z12").Value = [double]"20.75245245245293"
# Row 13
$ws.Range("F13").Value = [double]"25.8500000000006"
$ws.Range("H13").Value = [double]"0.001390610279246141"
$ws.Range("I13").Value = [double]"0.001390610279246141"
$ws.Range("L13").Value = [double]"34.76951870147207"
$ws.Range("M13").Value = "[12.582365563385764, 56.95667183955837]"
$ws.Range("N13").Value = [double]"0.002849327349497166"
$ws.Range("O13").Value = [double]"0.002849327349497166"
$ws.Range("P13").Value = [double]"1.66671081721981"
$ws.Range("Q13").Value = "[0.8365001460008852, 2.496921488438735]"
$ws.Range("R13").Value = [double]"0.0002037057853825353"
$ws.Range("S13").Value = [double]"0.0002037057853825353"
$ws.Range("T13").Value = [double]"61.73909555954534"
$ws.Range("U13").Value = "[48.56368301213537, 74.91450810695531]"
$ws.Range("V13").Value = [double]"3.08397751780376e-12"
$ws.Range("W13").Value = [double]"3.08397751780376e-12"
$ws.Range("X13").Value = [double]"18.99289289289333"
$ws.Range("Y13").Value = [double]"15.57727727727764"
$ws.Range("Z13").Value = [double]"22.40850850850903"
